$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad") for rows 2-9 from 45204 to 45207 (2023-10-05 -> 2023-10-08)
$ws.Range("C2:C9").Value = 45207
